$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns with updated values (F and J remain 0 and are unchanged)
$cols = @("B","C","D","E","G","H","I","K")

# New values for rows 2..25 (corresponding to data rows 0..23), one row per record,
# values in column order matching $cols above.
$data = @(
    @(0.3532548847276757,0.09797019725814948,0.078275756824425,0.4096749627931899,1.231866220413735,1.119606534948758,0.9479722294406159,0.4397157291636518),
    @(0.317814146708173,0.08626813890626295,0.07102748671820791,0.3573579241509464,1.198169773588859,1.108746642098907,0.9349257697296238,0.3930083870752128),
    @(0.2962143582678607,0.07911324692055643,0.06661546065564039,0.3253450692656088,1.17815168649895,1.102574353334774,0.9273496079584191,0.3645142541847406),
    @(0.2874525756803337,0.07620492130918421,0.06482709650134666,0.3123250101438799,1.170161772679378,1.100183275447108,0.9243708959740289,0.3529485775373189),
    @(0.2860001180399081,0.07572243286305991,0.06453071542171074,0.3101645110505444,1.168845140177865,1.099793725133253,0.9238828268853254,0.3510308660880241),
    @(0.296096030822099,0.07907399478443722,0.06659130351522435,0.3251693761003054,1.178043254739379,1.102541604264559,0.9273089969520001,0.3643580901006658),
    @(0.341001392680198,0.09392890199691806,0.07576851563985088,0.3916118538413684,1.220107498237923,1.11575886271163,0.9433832779957356,0.4235724861546259),
    @(0.4303507902402259,0.1233118304149627,0.09407525260088789,0.5228851543253228,1.307987707445079,1.145635176522944,0.9783814241279458,0.5411813646764756),
    @(0.4968077217168627,0.1450732468407239,0.107723042641183,0.6200862721951239,1.375935573942456,1.170034497764334,1.006258514832183,0.6285434326091206),
    @(0.5272231023531901,0.1550152803821447,0.1139767103588696,0.6645043050562549,1.407602314675557,1.181674416706812,1.019420291403662,0.6685048858417417),
    @(0.5387673999024969,0.1587865260777619,0.1163514456373349,0.6813558169657625,1.419704093723965,1.186160487618451,1.024474124982575,0.6836695892791909),
    @(0.5362799397754827,0.1579740311040894,0.1158397083747786,0.6777251071516019,1.41709283183144,1.185190841240683,1.023382578237133,0.6804021597361611),
    @(0.5281723231676381,0.1553254128839114,0.1141719482396439,0.6658900460023887,1.408595717086683,1.182041916277001,1.019834671429109,0.669751847318139),
    @(0.5232096517583216,0.1537039013318235,0.1131512605626881,0.6586448889827352,1.403405392185277,1.180123320622727,1.017670580466003,0.6632324287517122),
    @(0.494823722282888,0.1444243936981877,0.1073152695804964,0.6171877017639105,1.373881430806506,1.16928472625014,1.00540808568347,0.6259363298707967),
    @(0.4774571027479055,0.1387428389192564,0.1037467338188236,0.591808245056086,1.355964312872061,1.162774440034013,0.9980089535938106,0.6031131040369928),
    @(0.4674855674950891,0.1354789568699459,0.1016984563676573,0.5772293974307416,1.34573001882444,1.15908074060431,0.9937983400951538,0.5900064638586855),
    @(0.4641123390586301,0.1343745407108088,0.1010056715894621,0.5722963940077364,1.342277040912535,1.15783883327137,0.9923804363188111,0.5855723158415742),
    @(0.4793040185253119,0.1393472343349629,0.1041261699360803,0.5945079735227807,1.357864245767047,1.163462204869347,0.9987919238556913,0.6055405317334532),
    @(0.5305530014961448,0.1561031999507065,0.1146616296651501,0.6693654189708269,1.411088521472124,1.182964702733585,1.020874879087543,0.6728792248893285),
    @(0.5642026410085066,0.1670917072460725,0.1215857221585708,0.7184731026389528,1.446517082146698,1.196167287807555,1.035714269632138,0.7170767134998925),
    @(0.5462288927553232,0.1612234078675954,0.1178866404372769,0.6922457352216469,1.42754884420043,1.189078857307521,1.027756749372799,0.6934703222958092),
    @(0.4784689881579425,0.1390739793427542,0.1039546165320786,0.5932873886833647,1.357005079083734,1.163151113320652,0.9984378083396024,0.6044430458142074),
    @(0.4060387063375401,0.1153339024261015,0.08908868760214261,0.4872512108890419,1.283626722099712,1.13712541292935,0.9685367247694856,0.5092007997978101)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = 2 + $i
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $ws.Range($cols[$j] + $rowNum).Value = $rowVals[$j]
    }
}
